$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header for new column F, copying the formatting from the E1 header cell
$ws.Cells.Item(1, 5).Copy()
$ws.Cells.Item(1, 6).PasteSpecial(-4122)
$ws.Cells.Item(1, 6).Value = "time_taken"

$timestamps = @(
    "2021-10-05 13:41:33.046171",
    "2021-10-05 13:41:33.046181",
    "2021-10-05 13:41:33.046185",
    "2021-10-05 13:41:33.046187",
    "2021-10-05 13:41:33.046190",
    "2021-10-05 13:41:33.046192",
    "2021-10-05 13:41:33.046195",
    "2021-10-05 13:41:33.046197",
    "2021-10-05 13:41:33.046200",
    "2021-10-05 13:41:33.046202",
    "2021-10-05 13:41:33.046205",
    "2021-10-05 13:41:33.046207",
    "2021-10-05 13:41:33.046210",
    "2021-10-05 13:41:33.046212",
    "2021-10-05 13:41:33.046215",
    "2021-10-05 13:41:33.046217",
    "2021-10-05 13:41:33.046220",
    "2021-10-05 13:41:33.046223",
    "2021-10-05 13:41:33.046225",
    "2021-10-05 13:41:33.046228",
    "2021-10-05 13:41:33.046231",
    "2021-10-05 13:41:33.046233",
    "2021-10-05 13:41:33.046235",
    "2021-10-05 13:41:33.046238",
    "2021-10-05 13:41:33.046241",
    "2021-10-05 13:41:33.046243",
    "2021-10-05 13:41:33.046246",
    "2021-10-05 13:41:33.046248",
    "2021-10-05 13:41:33.046251",
    "2021-10-05 13:41:33.046253"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
